# 🚌 141: 31/12 10:46 LP1912+6203+6173
# Adds the latest scrape batch (run at 07:46:38 / 07:46:49 local time) to the
# three tracking sheets and refreshes the "last updated" / "total rows"
# banner cells on each sheet.

$wb = $excel.ActiveWorkbook

$newUpdated = "Última actualización: 31/12/2025 07:46:49"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912 — 14 new rows (717-730), banner text updated
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = $newUpdated
$ws1.Cells.Item(3, 1).Value = "Total filas: 729"

$sheet1Rows = @(
    @("07:46:38", "07:51", "15_ABASTO",      5,  "LP1912", "31/12/2025"),
    @("07:46:38", "07:58", "23_HERNANDEZ",   12, "LP1912", "31/12/2025"),
    @("07:46:38", "08:01", "16_SANTA ANA",   15, "LP1912", "31/12/2025"),
    @("07:46:38", "08:03", "17_ROMERO",      17, "LP1912", "31/12/2025"),
    @("07:46:38", "08:11", "16_SANTA ANA",   25, "LP1912", "31/12/2025"),
    @("07:46:38", "08:14", "10_OLMOS",       28, "LP1912", "31/12/2025"),
    @("07:46:38", "08:15", "17_ROMERO",      29, "LP1912", "31/12/2025"),
    @("07:46:38", "08:29", "14_ABASTO",      43, "LP1912", "31/12/2025"),
    @("07:46:38", "08:44", "10_OLMOS",       58, "LP1912", "31/12/2025"),
    @("07:46:38", "08:49", "16_SANTA ANA",   63, "LP1912", "31/12/2025"),
    @("07:46:38", "09:02", "17X38_ROMERO",   76, "LP1912", "31/12/2025"),
    @("07:46:38", "09:02", "23_HERNANDEZ",   76, "LP1912", "31/12/2025"),
    @("07:46:38", "09:14", "11_ETCHEVERRY",  88, "LP1912", "31/12/2025"),
    @("07:46:38", "09:16", "27_EL RETIRO",   90, "LP1912", "31/12/2025")
)

$r = 717
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = ""
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
    $ws1.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215 — banner text only (no new rows)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = $newUpdated

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173 — 1 new row (90), banner text updated
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = $newUpdated
$ws3.Cells.Item(3, 1).Value = "Total filas: 89"

$ws3.Cells.Item(90, 1).Value = ""
$ws3.Cells.Item(90, 2).Value = "31/12/2025"
$ws3.Cells.Item(90, 3).Value = "07:46:49"
$ws3.Cells.Item(90, 4).Value = "08:09"
$ws3.Cells.Item(90, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(90, 6).Value = 23
$ws3.Cells.Item(90, 7).Value = "L6173"
